# Reverts an erroneous merge:
#   - removes the spurious "Maize (home consumed)" row that had been
#     merged onto sheet "Y" (B2/C2), restoring the row-index sequence
#     1..5 down column A
#   - restores sheet selections / the active sheet & tab as they were
#     before the merge (sheet "A" becomes the active tab again)

$wb = $excel.ActiveWorkbook

$wsY  = $wb.Worksheets.Item("Y")
$wsA  = $wb.Worksheets.Item("A")
$wsVA = $wb.Worksheets.Item("VA")

# --- Sheet "Y": drop the merged-in Maize (home consumed) entry and
#     fill the row index column back in for rows 3-6 ---
$wsY.Range("B2:C2").Clear()
$wsY.Cells.Item(3, 1).Value = 2
$wsY.Cells.Item(4, 1).Value = 3
$wsY.Cells.Item(5, 1).Value = 4
$wsY.Cells.Item(6, 1).Value = 5

# --- Restore each sheet's selection ---
$wsY.Range("A4").Select()
$wsVA.Range("D21").Select()

# --- Sheet "A" is the active tab again (also flips tabSelected on its
#     sheetView and clears it from whichever sheet previously had it) ---
$wsA.Activate()
$wsA.Range("C2").Select()

# Scroll sheet "A" so row 25 is at the top of the view (best effort —
# harmless no-op if the host doesn't persist window scroll position).
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
